# rolling-ta EMA fix:
# Column C was computed as a 14-period simple moving AVERAGE for every row.
# Replace it with the correct EMA recursion, seeded off the existing C14
# average: C(r) = (B(r) - C(r-1)) * (2/15) + C(r-1), for every row from 15
# down to 200. C14 itself (the seed average) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 15; $r -le 200; $r++) {
    $prev = $r - 1
    $ws.Range("C$r").Formula = "=((B$r-C$prev)*(2/15))+C$prev"
}

# Restore the cursor/selection to where the author left it after the edit.
$ws.Range("G14").Select() | Out-Null
